$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.092.89"
Set-TextValue $ws.Range("E2") "  -1.12%  "
Set-TextValue $ws.Range("D3") "2.605.75"
Set-TextValue $ws.Range("E3") "  -0.64%  "
Set-TextValue $ws.Range("E4") "  -0.04%  "
Set-TextValue $ws.Range("D5") "590.79"
Set-TextValue $ws.Range("E5") "  -2.40%  "
Set-TextValue $ws.Range("D6") "150.23"
Set-TextValue $ws.Range("E6") "  -2.78%  "
Set-TextValue $ws.Range("E7") "  +0.00%  "
Set-TextValue $ws.Range("D8") "0.547"
Set-TextValue $ws.Range("E8") "  -0.81%  "
Set-TextValue $ws.Range("D9") "2.605.55"
Set-TextValue $ws.Range("E9") "  -0.58%  "
Set-TextValue $ws.Range("D10") "0.128"
Set-TextValue $ws.Range("E10") "  -0.91%  "
Set-TextValue $ws.Range("E11") "  -0.04%  "
Set-TextValue $ws.Range("E12") "  -1.62%  "
Set-TextValue $ws.Range("D13") "0.344"
Set-TextValue $ws.Range("E13") "  -3.14%  "
Set-TextValue $ws.Range("D14") "27.35"
Set-TextValue $ws.Range("E14") "  -1.74%  "
Set-TextValue $ws.Range("D15") "3.070.63"
Set-TextValue $ws.Range("E15") "  -1.20%  "
Set-TextValue $ws.Range("D16") "0.0000181"
Set-TextValue $ws.Range("E16") "  -4.21%  "
Set-TextValue $ws.Range("D17") "66.909.73"
Set-TextValue $ws.Range("E17") "  -1.32%  "
Set-TextValue $ws.Range("D18") "2.599.68"
Set-TextValue $ws.Range("E18") "  -0.49%  "
Set-TextValue $ws.Range("D19") "364.53"
Set-TextValue $ws.Range("E19") "  -0.59%  "
Set-TextValue $ws.Range("D20") "11.05"
Set-TextValue $ws.Range("E20") "  -1.69%  "
Set-TextValue $ws.Range("D21") "7.34"
Set-TextValue $ws.Range("E21") "  -4.43%  "
Set-TextValue $ws.Range("E22") "  -0.35%  "
Set-TextValue $ws.Range("D23") "4.87"
Set-TextValue $ws.Range("E23") "  -2.08%  "
Set-TextValue $ws.Range("D24") "2.06"
Set-TextValue $ws.Range("E24") "  -0.59%  "
Set-TextValue $ws.Range("D25") "72.90"
Set-TextValue $ws.Range("E25") "  +3.51%  "
Set-TextValue $ws.Range("E26") "  +0.04%  "
Set-TextValue $ws.Range("D27") "9.96"
Set-TextValue $ws.Range("E27") "  +0.26%  "
Set-TextValue $ws.Range("B28") "Bittensor"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D28") "590.07"
Set-TextValue $ws.Range("E28") "  +1.61%  "
Set-TextValue $ws.Range("B29") "WrappedeETH"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D29") "2.733.17"
Set-TextValue $ws.Range("E29") "  -0.48%  "
Set-TextValue $ws.Range("E30") "  -0.09%  "
Set-TextValue $ws.Range("D31") "0.0₃0988"
Set-TextValue $ws.Range("E31") "  -6.10%  "
Set-TextValue $ws.Range("D32") "1.37"
Set-TextValue $ws.Range("E32") "  -5.12%  "
Set-TextValue $ws.Range("D33") "7.65"
Set-TextValue $ws.Range("E33") "  -3.47%  "
Set-TextValue $ws.Range("D34") "1.81"
Set-TextValue $ws.Range("E34") "  -3.50%  "
Set-TextValue $ws.Range("E35") "  -0.03%  "
Set-TextValue $ws.Range("E36") "  -4.90%  "
Set-TextValue $ws.Range("D37") "1.49"
Set-TextValue $ws.Range("E37") "  -2.67%  "
Set-TextValue $ws.Range("D38") "155.65"
Set-TextValue $ws.Range("E38") "  -1.26%  "
Set-TextValue $ws.Range("D39") "19.00"
Set-TextValue $ws.Range("E39") "  -2.37%  "
Set-TextValue $ws.Range("E40") "  -1.42%  "
Set-TextValue $ws.Range("E41") "  -0.46%  "
Set-TextValue $ws.Range("D42") "5.21"
Set-TextValue $ws.Range("E42") "  -3.18%  "
Set-TextValue $ws.Range("B43") "dogwifhat"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D43") "2.56"
Set-TextValue $ws.Range("E43") "  -3.31%  "
Set-TextValue $ws.Range("B44") "WhiteBITCoin"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D44") "17.08"
Set-TextValue $ws.Range("E44") "  +3.93%  "
Set-TextValue $ws.Range("E45") "  -0.06%  "
Set-TextValue $ws.Range("D46") "153.25"
Set-TextValue $ws.Range("E46") "  -2.63%  "
Set-TextValue $ws.Range("E47") "  -1.90%  "
Set-TextValue $ws.Range("D48") "3.71"
Set-TextValue $ws.Range("E48") "  -1.41%  "
Set-TextValue $ws.Range("E49") "  -2.75%  "
Set-TextValue $ws.Range("D50") "0.0779"
Set-TextValue $ws.Range("E50") "  -1.48%  "
Set-TextValue $ws.Range("D51") "21.52"
Set-TextValue $ws.Range("E51") "  +2.31%  "
